$wb = $excel.ActiveWorkbook

# Rename Sheet1 to Table_Names
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Table_Names"

# Clear existing content and set new values
$ws1.Cells.Clear()
$ws1.Range("A1").Value = "T1"
$ws1.Range("A2").Value = "T2"

# Add a new sheet named Field_Names, placed after Table_Names
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Field_Names"

$ws1.Select()
$ws1.Range("D3").Select()
